$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '68.333.24'
$ws.Range("E2").Value = '  -1.82%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.836.16'
$ws.Range("E3").Value = '  -1.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '601.53'
$ws.Range("E5").Value = '  -0.67%  '

# Row 6
Set-TextValue $ws.Range("D6") '169.77'
$ws.Range("E6").Value = '  -0.09%  '

# Row 7
Set-TextValue $ws.Range("D7") '3.837.29'
$ws.Range("E7").Value = '  -1.01%  '

# Row 9
$ws.Range("E9").Value = '  -1.35%  '

# Row 10
$ws.Range("E10").Value = '  -2.21%  '

# Row 11
$ws.Range("E11").Value = '  +1.51%  '

# Row 12
$ws.Range("E12").Value = '  -1.99%  '

# Row 13
$ws.Range("E13").Value = '  +3.19%  '

# Row 14
Set-TextValue $ws.Range("D14") '37.25'
$ws.Range("E14").Value = '  -2.65%  '

# Row 15
Set-TextValue $ws.Range("D15") '4.479.57'
$ws.Range("E15").Value = '  -1.41%  '

# Row 16
Set-TextValue $ws.Range("D16") '3.835.87'
$ws.Range("E16").Value = '  -1.73%  '

# Row 17
Set-TextValue $ws.Range("D17") '68.427.68'
$ws.Range("E17").Value = '  -1.73%  '

# Row 18
Set-TextValue $ws.Range("D18") '18.50'
$ws.Range("E18").Value = '  -0.82%  '

# Row 19
Set-TextValue $ws.Range("D19") '7.42'
$ws.Range("E19").Value = '  -2.99%  '

# Row 20
$ws.Range("E20").Value = '  -0.84%  '

# Row 21
Set-TextValue $ws.Range("D21") '11.09'
$ws.Range("E21").Value = '  -0.26%  '

# Row 22
Set-TextValue $ws.Range("D22") '469.50'
$ws.Range("E22").Value = '  -4.25%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.738'
$ws.Range("E23").Value = '  -1.53%  '

# Row 24
$ws.Range("E24").Value = '  -3.88%  '

# Row 25
Set-TextValue $ws.Range("D25") '83.24'
$ws.Range("E25").Value = '  -2.51%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.24'
$ws.Range("E26").Value = '  -2.65%  '

# Row 27
$ws.Range("E27").Value = '  -0.51%  '

# Row 28
Set-TextValue $ws.Range("D28") '10.06'
$ws.Range("E28").Value = '  -0.81%  '

# Row 29
$ws.Range("E29").Value = '  -0.10%  '

# Row 30
$ws.Range("E30").Value = '  -0.09%  '

# Row 31
Set-TextValue $ws.Range("D31") '3.986.03'
$ws.Range("E31").Value = '  -1.31%  '

# Row 32
$ws.Range("E32").Value = '  -1.40%  '

# Row 33
Set-TextValue $ws.Range("D33") '31.66'
$ws.Range("E33").Value = '  -0.86%  '

# Row 34
Set-TextValue $ws.Range("D34") '2.32'
$ws.Range("E34").Value = '  -4.13%  '

# Row 35
Set-TextValue $ws.Range("D35") '9.48'
$ws.Range("E35").Value = '  -1.26%  '

# Row 36
Set-TextValue $ws.Range("D36") '3.801.89'
$ws.Range("E36").Value = '  -1.29%  '

# Row 37
$ws.Range("E37").Value = '  -1.79%  '

# Row 38
Set-TextValue $ws.Range("D38") '3.71'
$ws.Range("E38").Value = '  +11.57%  '

# Row 39
$ws.Range("E39").Value = '  -0.70%  '

# Row 40
$ws.Range("E40").Value = '  -1.62%  '

# Row 41
$ws.Range("E41").Value = '  -2.68%  '

# Row 42
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.316'
$ws.Range("E43").Value = '  -3.64%  '

# Row 44
$ws.Range("E44").Value = '  -5.41%  '

# Row 45
$ws.Range("E45").Value = '  +0.87%  '

# Row 46
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue $ws.Range("D46") '0.000295'
$ws.Range("E46").Value = '  +7.55%  '

# Row 47
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D47") '418.98'
$ws.Range("E47").Value = '  -3.98%  '

# Row 49
$ws.Range("E49").Value = '  -2.13%  '

# Row 50
$ws.Range("E50").Value = '  -2.36%  '

# Row 51
Set-TextValue $ws.Range("D51") '141.75'
$ws.Range("E51").Value = '  -1.52%  '
